$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
# This shifts the old N/O/P ("Late", "Outstanding", "Disbursement") columns one to the
# right (-> O/P/Q) and leaves a new blank column N in between "In Advance" (M) and
# "Late" (now O), matching column M's width/look.
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and move the selection, as recorded
# in the workbook when it was last saved.
$null = $wsSchedule.Activate()
$null = $wsSchedule.Range("K17").Select()
